$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 164; this shifts the existing rows 164-252 down to 165-253
$ws.Rows(164).Insert()

# Populate the newly inserted row 164 with the new data record
$ws.Cells.Item(164, 1).Value  = 10
$ws.Cells.Item(164, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(164, 3).Value  = "La Araucanía"
$ws.Cells.Item(164, 4).Value  = 44572
$ws.Cells.Item(164, 5).Value  = 9
$ws.Cells.Item(164, 6).Value  = 100112009
$ws.Cells.Item(164, 7).Value  = "Acelga"
$ws.Cells.Item(164, 8).Value  = "Sin especificar"
$ws.Cells.Item(164, 9).Value  = "Primera"
$ws.Cells.Item(164, 10).Value = 65
$ws.Cells.Item(164, 11).Value = 8000
$ws.Cells.Item(164, 12).Value = 8000
$ws.Cells.Item(164, 13).Value = 8000
$ws.Cells.Item(164, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(164, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(164, 16).Value = 667
$ws.Cells.Item(164, 17).Value = 12
$ws.Cells.Item(164, 18).Value = "Hortaliza"
